$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 new values
$ws.Range("A12").Value = 112181852
$ws.Range("B12").Value = 77671
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 185
$ws.Range("F12").Value = "Violettgrå tagellav"
$ws.Range("G12").Value = "Bryoria nadvornikiana"
$ws.Range("H12").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q12").Value = 374954
$ws.Range("R12").Value = 6870892

# Row 13 new values
$ws.Range("A13").Value = 112182046
$ws.Range("B13").Value = 90812
$ws.Range("E13").Value = 4366
$ws.Range("F13").Value = "Skarp dropptaggsvamp"
$ws.Range("G13").Value = "Hydnellum peckii"
$ws.Range("H13").Value = "Banker"
$ws.Range("Q13").Value = 374850
$ws.Range("R13").Value = 6871061

# Row 14 new values
$ws.Range("A14").Value = 112182730
$ws.Range("B14").Value = 95693
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 221941
$ws.Range("F14").Value = "Plattlummer"
$ws.Range("G14").Value = "Lycopodium complanatum"
$ws.Range("H14").Value = "L."
$ws.Range("Q14").Value = 375047
$ws.Range("R14").Value = 6871264
